$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: Total 13 -> 23, Book Code B-000024 -> B-000001
$ws.Range("D2").Value = 23.0
$ws.Range("C2").Value = "B-000001"

# Add new row 3 with values mirroring row 2's warehouse/supplier/address info
$ws.Range("B3").Value = "W-000001"
$ws.Range("C3").Value = "B-000001"
$ws.Range("D3").Value = 0.0
$ws.Range("E3").Value = 0.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = "ORI"
$ws.Range("H3").Value = "Sohel1"
$ws.Range("I3").Value = "Hajipara"
$ws.Range("J3").Value = "Hajipara 2"

# Match row 3's formatting to row 2 (same style applied across the row)
$ws.Range("B2:J2").Copy()
$ws.Range("B3:J3").PasteSpecial(-4122)
